$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2028985507246377
$ws.Range("C2").Value = 0.5594202898550724
$ws.Range("J2").Value = 0.02318840579710145
$ws.Range("P2").Value = 0.1507246376811594
$ws.Range("S2").Value = 0.06376811594202898
$ws.Range("B3").Value = 0.004901960784313725
$ws.Range("C3").Value = 0.05882352941176471
$ws.Range("J3").Value = 0.03431372549019608
$ws.Range("P3").Value = 0.7450980392156863
$ws.Range("S3").Value = 0.1568627450980392
$ws.Range("J4").Value = 0.01886792452830189
$ws.Range("P4").Value = 0.8301886792452831
$ws.Range("S4").Value = 0.1509433962264151
$ws.Range("B6").Value = 0.05855855855855856
$ws.Range("D6").Value = 0.02252252252252252
$ws.Range("E6").Value = 0.004504504504504504
$ws.Range("F6").Value = 0.06306306306306306
$ws.Range("J6").Value = 0.3018018018018018
$ws.Range("O6").Value = 0.004504504504504504
$ws.Range("Q6").Value = 0.1576576576576577
$ws.Range("R6").Value = 0.04504504504504504
$ws.Range("S6").Value = 0.3423423423423423
$ws.Range("B7").Value = 0.1187739463601533
$ws.Range("D7").Value = 0.01149425287356322
$ws.Range("E7").Value = 0.003831417624521073
$ws.Range("F7").Value = 0.04980842911877394
$ws.Range("J7").Value = 0.1264367816091954
$ws.Range("O7").Value = 0.01532567049808429
$ws.Range("Q7").Value = 0.1532567049808429
$ws.Range("R7").Value = 0.08812260536398467
$ws.Range("S7").Value = 0.4329501915708812
$ws.Range("B8").Value = 0.1044776119402985
$ws.Range("D8").Value = 0.02132196162046908
$ws.Range("E8").Value = 0.002132196162046908
$ws.Range("F8").Value = 0.04264392324093817
$ws.Range("J8").Value = 0.1002132196162047
$ws.Range("O8").Value = 0.01492537313432836
$ws.Range("Q8").Value = 0.1556503198294243
$ws.Range("R8").Value = 0.1023454157782516
$ws.Range("S8").Value = 0.4562899786780384
$ws.Range("B9").Value = 0.09836065573770492
$ws.Range("D9").Value = 0.00546448087431694
$ws.Range("F9").Value = 0.0546448087431694
$ws.Range("J9").Value = 0.1038251366120219
$ws.Range("O9").Value = 0.03278688524590164
$ws.Range("Q9").Value = 0.180327868852459
$ws.Range("R9").Value = 0.06010928961748634
$ws.Range("S9").Value = 0.4644808743169399
$ws.Range("B10").Value = 0.1173541963015647
$ws.Range("D10").Value = 0.02489331436699858
$ws.Range("E10").Value = 0.0007112375533428165
$ws.Range("F10").Value = 0.05974395448079658
$ws.Range("J10").Value = 0.1130867709815078
$ws.Range("O10").Value = 0.01564722617354196
$ws.Range("Q10").Value = 0.2226173541963016
$ws.Range("R10").Value = 0.07041251778093884
$ws.Range("S10").Value = 0.3755334281650071
$ws.Range("G11").Value = 0.1517412935323383
$ws.Range("J11").Value = 0.07462686567164178
$ws.Range("K11").Value = 0.1965174129353234
$ws.Range("L11").Value = 0.5597014925373134
$ws.Range("S11").Value = 0.01741293532338309
$ws.Range("G12").Value = 0.7627118644067796
$ws.Range("J12").Value = 0.173728813559322
$ws.Range("K12").Value = 0.00423728813559322
$ws.Range("L12").Value = 0.0423728813559322
$ws.Range("S12").Value = 0.01694915254237288
$ws.Range("G13").Value = 0.5625
$ws.Range("J13").Value = 0.3958333333333333
$ws.Range("S13").Value = 0.04166666666666666
$ws.Range("F15").Value = 0.03292181069958848
$ws.Range("H15").Value = 0.1316872427983539
$ws.Range("I15").Value = 0.04938271604938271
$ws.Range("K15").Value = 0.06584362139917696
$ws.Range("M15").Value = 0.01234567901234568
$ws.Range("O15").Value = 0.06584362139917696
$ws.Range("S15").Value = 0.2345679012345679
$ws.Range("F16").Value = 0.01224489795918367
$ws.Range("H16").Value = 0.1591836734693877
$ws.Range("I16").Value = 0.09387755102040816
$ws.Range("J16").Value = 0.3959183673469387
$ws.Range("K16").Value = 0.1265306122448979
$ws.Range("M16").Value = 0.01224489795918367
$ws.Range("O16").Value = 0.07346938775510205
$ws.Range("S16").Value = 0.1265306122448979
$ws.Range("F17").Value = 0.02626262626262626
$ws.Range("H17").Value = 0.1474747474747475
$ws.Range("I17").Value = 0.08686868686868687
$ws.Range("J17").Value = 0.4141414141414141
$ws.Range("K17").Value = 0.1232323232323232
$ws.Range("M17").Value = 0.02424242424242424
$ws.Range("O17").Value = 0.06060606060606061
$ws.Range("S17").Value = 0.1171717171717172
$ws.Range("F18").Value = 0.01595744680851064
$ws.Range("H18").Value = 0.1223404255319149
$ws.Range("I18").Value = 0.05319148936170213
$ws.Range("J18").Value = 0.4361702127659575
$ws.Range("K18").Value = 0.1117021276595745
$ws.Range("M18").Value = 0.02127659574468085
$ws.Range("O18").Value = 0.0797872340425532
$ws.Range("S18").Value = 0.1595744680851064
$ws.Range("F19").Value = 0.02036516853932584
$ws.Range("H19").Value = 0.2127808988764045
$ws.Range("I19").Value = 0.07092696629213484
$ws.Range("J19").Value = 0.3546348314606741
$ws.Range("K19").Value = 0.1320224719101123
$ws.Range("M19").Value = 0.01896067415730337
$ws.Range("N19").Value = 0.0007022471910112359
$ws.Range("O19").Value = 0.06601123595505617
$ws.Range("S19").Value = 0.1235955056179775
